$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("JMX Files")

$ws.Range("C2").Value = "Dog API Test.jmx"
$ws.Range("C3").Value = "Dog API Test.jmx"

$ws.Range("C3").Select()
